# The document body is a single paragraph whose run contains a sequence of
# <w:t> text chunks separated by manual line breaks (<w:br/>), including one
# trailing break after the final text chunk. The edit collapses this into a
# single <w:t> run: every internal line break becomes a plain space, and the
# final (trailing) line break is simply dropped.

$d = $word.ActiveDocument

# Drop the trailing manual line break (the very last character before the
# paragraph mark) - it has no replacement in the target text.
$full = $d.Content.Text
$lastBreak = $d.Range($full.Length - 2, $full.Length - 1)
if ($lastBreak.Text -eq [string][char]11) {
    $lastBreak.Delete()
}

# Replace every remaining manual line break with a single space, merging all
# the text chunks into one run / one <w:t>.
$d.Content.Find.Execute("^l", $false, $false, $false, $false, $false, $true, 1, $false, " ", 2) | Out-Null
